$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 14531.5
$ws.Cells.Item(2, 9).Value = 10000
$ws.Cells.Item(2, 10).Value = 15178.857
$ws.Cells.Item(2, 11).Value = 10000
$ws.Cells.Item(2, 12).Value = 15178.857
$ws.Cells.Item(2, 13).Value = -9887
$ws.Cells.Item(2, 14).Value = -15404.857
$ws.Cells.Item(6, 8).Value = 161.14285
$ws.Cells.Item(6, 9).Value = 161.14285
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 483.42855
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = -371.42855
$ws.Cells.Item(6, 14).ClearContents()
$ws.Cells.Item(8, 8).Value = 1198.5834
$ws.Cells.Item(8, 9).Value = 1198.5834
$ws.Cells.Item(8, 11).Value = 3595.7502
$ws.Cells.Item(8, 13).Value = -3456.7502
$ws.Cells.Item(46, 8).Value = 111115090
$ws.Cells.Item(46, 10).Value = 111115090
$ws.Cells.Item(46, 12).Value = 333345270
$ws.Cells.Item(46, 14).Value = -333345508
$ws.Cells.Item(52, 8).Value = 37038140
$ws.Cells.Item(52, 9).Value = 700
$ws.Cells.Item(52, 10).Value = 47620264
$ws.Cells.Item(52, 11).Value = 2100
$ws.Cells.Item(52, 12).Value = 142860792
$ws.Cells.Item(52, 13).Value = -1940
$ws.Cells.Item(52, 14).Value = -142861112
$ws.Cells.Item(60, 8).Value = 111115090
$ws.Cells.Item(60, 10).Value = 111115090
$ws.Cells.Item(60, 12).Value = 333345270
$ws.Cells.Item(60, 14).Value = -333346238
$ws.Cells.Item(86, 8).Value = 2885.2144
$ws.Cells.Item(86, 9).Value = 2319.5
$ws.Cells.Item(86, 11).Value = 2319.5
$ws.Cells.Item(86, 13).Value = -1196.5
$ws.Cells.Item(89, 8).Value = 2885.2144
$ws.Cells.Item(89, 9).Value = 2319.5
$ws.Cells.Item(89, 11).Value = 11597.5
$ws.Cells.Item(89, 13).Value = -5981.5
$ws.Cells.Item(106, 8).Value = 4677.5713
$ws.Cells.Item(106, 9).Value = 2123.8333
$ws.Cells.Item(106, 11).Value = 2123.8333
$ws.Cells.Item(106, 13).Value = -1492.8333
$ws.Cells.Item(111, 8).Value = 11113423
$ws.Cells.Item(111, 9).Value = 3179.75
$ws.Cells.Item(111, 10).Value = 20001618
$ws.Cells.Item(111, 11).Value = 9539.25
$ws.Cells.Item(111, 12).Value = 60004854
$ws.Cells.Item(111, 13).Value = -6472.25
$ws.Cells.Item(111, 14).Value = -60010988
$ws.Cells.Item(116, 8).Value = 8310.799999999999
$ws.Cells.Item(116, 9).Value = 7185
$ws.Cells.Item(116, 11).Value = 7185
$ws.Cells.Item(116, 13).Value = -3743

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 558972.5
$ws.Cells.Item(74, 9).Value = 1111770.5
$ws.Cells.Item(74, 11).Value = 1111770.5
$ws.Cells.Item(74, 13).Value = -1110896.5
$ws.Cells.Item(77, 8).Value = 558972.5
$ws.Cells.Item(77, 9).Value = 1111770.5
$ws.Cells.Item(77, 11).Value = 5558852.5
$ws.Cells.Item(77, 13).Value = -5554484.5
$ws.Cells.Item(97, 8).Value = 1545970.9
$ws.Cells.Item(97, 9).Value = 1686440.9
$ws.Cells.Item(97, 11).Value = 1686440.9
$ws.Cells.Item(97, 13).Value = -1685944.9
$ws.Cells.Item(122, 8).Value = 1534.1852
$ws.Cells.Item(122, 9).Value = 1062.3478
$ws.Cells.Item(122, 11).Value = 3187.0434
$ws.Cells.Item(122, 13).Value = -737.0434

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 2166.8333
$ws.Cells.Item(22, 9).Value = 2400.2
$ws.Cells.Item(22, 10).Value = 1000
$ws.Cells.Item(22, 11).Value = 2400.2
$ws.Cells.Item(22, 12).Value = 1000
$ws.Cells.Item(22, 13).Value = -2227.2
$ws.Cells.Item(22, 14).Value = -1346
$ws.Cells.Item(74, 8).Value = 29900
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 13).ClearContents()
$ws.Cells.Item(77, 8).Value = 29900
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 13).ClearContents()
$ws.Cells.Item(94, 8).Value = 1887.3334
$ws.Cells.Item(94, 9).Value = 1369.5
$ws.Cells.Item(94, 10).Value = 2923
$ws.Cells.Item(94, 11).Value = 1369.5
$ws.Cells.Item(94, 12).Value = 2923
$ws.Cells.Item(94, 13).Value = -918.5
$ws.Cells.Item(94, 14).Value = -3825
$ws.Cells.Item(132, 8).Value = 98874.5
$ws.Cells.Item(132, 10).Value = 98874.5
$ws.Cells.Item(132, 12).Value = 98874.5
$ws.Cells.Item(132, 14).Value = -108994.5
$ws.Cells.Item(134, 8).Value = 2982.1035
$ws.Cells.Item(134, 9).Value = 1216.5652
$ws.Cells.Item(134, 10).Value = 9750
$ws.Cells.Item(134, 11).Value = 3649.6956
$ws.Cells.Item(134, 12).Value = 29250
$ws.Cells.Item(134, 13).Value = -1114.6956
$ws.Cells.Item(134, 14).Value = -34320

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).ClearContents()
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).ClearContents()
$ws.Cells.Item(99, 8).Value = 3840.125
$ws.Cells.Item(99, 9).Value = 3820.8572
$ws.Cells.Item(99, 11).Value = 3820.8572
$ws.Cells.Item(99, 13).Value = -2322.8572
$ws.Cells.Item(107, 8).Value = 1066.25
$ws.Cells.Item(107, 9).Value = 828.75
$ws.Cells.Item(107, 11).Value = 828.75
$ws.Cells.Item(107, 13).Value = 1091.25
$ws.Cells.Item(126, 8).Value = 3840.125
$ws.Cells.Item(126, 9).Value = 3820.8572
$ws.Cells.Item(126, 11).Value = 11462.5716
$ws.Cells.Item(126, 13).Value = -8992.571599999999
$ws.Cells.Item(132, 8).Value = 154349.12
$ws.Cells.Item(132, 9).Value = 3048.8333
$ws.Cells.Item(132, 11).Value = 9146.499899999999
$ws.Cells.Item(132, 13).Value = -6616.499899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10, 8).Value = 89.36364
$ws.Cells.Item(10, 9).Value = 81.333336
$ws.Cells.Item(10, 10).Value = 125.5
$ws.Cells.Item(10, 11).Value = 244.000008
$ws.Cells.Item(10, 12).Value = 376.5
$ws.Cells.Item(10, 13).Value = -105.000008
$ws.Cells.Item(10, 14).Value = -654.5
$ws.Cells.Item(19, 8).Value = 490
$ws.Cells.Item(19, 9).Value = 490
$ws.Cells.Item(19, 10).Value = 490
$ws.Cells.Item(19, 11).Value = 1470
$ws.Cells.Item(19, 12).Value = 1470
$ws.Cells.Item(19, 13).Value = -1296
$ws.Cells.Item(19, 14).Value = -1818
$ws.Cells.Item(23, 8).Value = 163.75
$ws.Cells.Item(23, 10).Value = 264
$ws.Cells.Item(23, 12).Value = 792
$ws.Cells.Item(23, 14).Value = -1262
$ws.Cells.Item(25, 8).Value = 765
$ws.Cells.Item(25, 9).Value = 20
$ws.Cells.Item(25, 10).Value = 1510
$ws.Cells.Item(25, 11).Value = 60
$ws.Cells.Item(25, 12).Value = 4530
$ws.Cells.Item(25, 13).Value = 109
$ws.Cells.Item(25, 14).Value = -4868
$ws.Cells.Item(29, 8).Value = 1700
$ws.Cells.Item(29, 9).Value = 50
$ws.Cells.Item(29, 11).Value = 150
$ws.Cells.Item(29, 13).Value = 127
$ws.Cells.Item(30, 8).Value = 765
$ws.Cells.Item(30, 9).Value = 20
$ws.Cells.Item(30, 10).Value = 1510
$ws.Cells.Item(30, 11).Value = 60
$ws.Cells.Item(30, 12).Value = 4530
$ws.Cells.Item(30, 13).Value = 42
$ws.Cells.Item(30, 14).Value = -4734
$ws.Cells.Item(62, 8).Value = 1738.9565
$ws.Cells.Item(65, 8).Value = 1738.9565
$ws.Cells.Item(128, 8).Value = 427432.5
$ws.Cells.Item(128, 9).Value = 427432.5
$ws.Cells.Item(128, 11).Value = 1282297.5
$ws.Cells.Item(128, 13).Value = -1277317.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1421.7059
$ws.Cells.Item(102, 9).Value = 1053
$ws.Cells.Item(102, 11).Value = 1053
$ws.Cells.Item(102, 13).Value = 569

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3366.9707
$ws.Cells.Item(22, 9).Value = 2184.2778
$ws.Cells.Item(22, 10).Value = 4697.5
$ws.Cells.Item(22, 11).Value = 2184.2778
$ws.Cells.Item(22, 12).Value = 4697.5
$ws.Cells.Item(22, 13).Value = -1889.2778
$ws.Cells.Item(22, 14).Value = -5287.5
$ws.Cells.Item(27, 8).Value = 3366.9707
$ws.Cells.Item(27, 9).Value = 2184.2778
$ws.Cells.Item(27, 10).Value = 4697.5
$ws.Cells.Item(27, 11).Value = 2184.2778
$ws.Cells.Item(27, 12).Value = 4697.5
$ws.Cells.Item(27, 13).Value = -2077.2778
$ws.Cells.Item(27, 14).Value = -4911.5
$ws.Cells.Item(122, 8).Value = 3558.8667
$ws.Cells.Item(122, 9).Value = 3491
$ws.Cells.Item(122, 11).Value = 10473
$ws.Cells.Item(122, 13).Value = -8023

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(32, 8).Value = 7750
$ws.Cells.Item(32, 9).Value = 7750
$ws.Cells.Item(32, 11).Value = 7750
$ws.Cells.Item(32, 13).Value = -7433
$ws.Cells.Item(81, 8).Value = 6478.3
$ws.Cells.Item(81, 9).Value = 2826.1428
$ws.Cells.Item(81, 11).Value = 5652.2856
$ws.Cells.Item(81, 13).Value = -4591.2856
$ws.Cells.Item(84, 8).Value = 6478.3
$ws.Cells.Item(84, 9).Value = 2826.1428
$ws.Cells.Item(84, 11).Value = 28261.428
$ws.Cells.Item(84, 13).Value = -22957.428
